$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (sending cluster, ligand, receptor, target cluster, plus stats)
# Row 2: Inflammatory-Mac -> ECs
$ws.Range("A2").Value = "Inflammatory-Mac"
$ws.Range("B2").Value = "Fgf15"
$ws.Range("C2").Value = "Fgfr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.018508
$ws.Range("H2").Value = 0.055524
$ws.Range("I2").Value = 0.2347866901774728
$ws.Range("J2").Value = 0.3151800006811757
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.330840333333333
$ws.Range("N2").Value = 6.992521
$ws.Range("O2").Value = 0.6715345129768794
$ws.Range("P2").Value = 0.7003397275969581
$ws.Range("Q2").Value = 0.04313919288933333
$ws.Range("R2").Value = 0.388252736004
$ws.Range("S2").Value = 0.1576673656417827
$ws.Range("T2").Value = 0.2207330758210636

# Row 3: Inflammatory-Mac -> FAPs
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("B3").Value = "Fgf15"
$ws.Range("C3").Value = "Fgfr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.018508
$ws.Range("H3").Value = 0.055524
$ws.Range("I3").Value = 0.2347866901774728
$ws.Range("J3").Value = 0.3151800006811757
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.657666
$ws.Range("N3").Value = 1.972998
$ws.Range("O3").Value = 0.1894790521235985
$ws.Range("P3").Value = 0.1976066831789769
$ws.Range("Q3").Value = 0.012172082328
$ws.Range("R3").Value = 0.109548740952
$ws.Range("S3").Value = 0.04448715950606454
$ws.Range("T3").Value = 0.06228167453895482

# Row 4: Inflammatory-Mac -> MuSCs
$ws.Range("A4").Value = "Inflammatory-Mac"
$ws.Range("B4").Value = "Fgf15"
$ws.Range("C4").Value = "Fgfr3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.018508
$ws.Range("H4").Value = 0.055524
$ws.Range("I4").Value = 0.2347866901774728
$ws.Range("J4").Value = 0.3151800006811757
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.42828
$ws.Range("N4").Value = 0.85656
$ws.Range("O4").Value = 0.123391035029171
$ws.Range("P4").Value = 0.08578923067523865
$ws.Range("Q4").Value = 0.00792660624
$ws.Range("R4").Value = 0.04755963744
$ws.Range("S4").Value = 0.02897057271207166
$ws.Range("T4").Value = 0.02703904978265926

# Row 5: Inflammatory-Mac -> Resolving-Mac
$ws.Range("A5").Value = "Inflammatory-Mac"
$ws.Range("B5").Value = "Fgf15"
$ws.Range("C5").Value = "Fgfr3"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.018508
$ws.Range("H5").Value = 0.055524
$ws.Range("I5").Value = 0.2347866901774728
$ws.Range("J5").Value = 0.3151800006811757
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.05413033333333334
$ws.Range("N5").Value = 0.162391
$ws.Range("O5").Value = 0.01559539987035126
$ws.Range("P5").Value = 0.01626435854882633
$ws.Range("Q5").Value = 0.001001844209333333
$ws.Range("R5").Value = 0.009016597883999999
$ws.Range("S5").Value = 0.003661592317553959
$ws.Range("T5").Value = 0.005126200538497968

# Row 6: MuSCs -> ECs
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Fgf15"
$ws.Range("C6").Value = "Fgfr3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.060321
$ws.Range("H6").Value = 0.120642
$ws.Range("I6").Value = 0.7652133098225272
$ws.Range("J6").Value = 0.6848199993188243
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.330840333333333
$ws.Range("N6").Value = 6.992521
$ws.Range("O6").Value = 0.6715345129768794
$ws.Range("P6").Value = 0.7003397275969581
$ws.Range("Q6").Value = 0.140598619747
$ws.Range("R6").Value = 0.843591718482
$ws.Range("S6").Value = 0.5138671473350968
$ws.Range("T6").Value = 0.4796066517758944

# Row 7: MuSCs -> FAPs
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Fgf15"
$ws.Range("C7").Value = "Fgfr3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.060321
$ws.Range("H7").Value = 0.120642
$ws.Range("I7").Value = 0.7652133098225272
$ws.Range("J7").Value = 0.6848199993188243
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.657666
$ws.Range("N7").Value = 1.972998
$ws.Range("O7").Value = 0.1894790521235985
$ws.Range("P7").Value = 0.1976066831789769
$ws.Range("Q7").Value = 0.039671070786
$ws.Range("R7").Value = 0.238026424716
$ws.Range("S7").Value = 0.144991892617534
$ws.Range("T7").Value = 0.1353250086400221

# Row 8: MuSCs -> MuSCs (new row)
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Fgf15"
$ws.Range("C8").Value = "Fgfr3"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.060321
$ws.Range("H8").Value = 0.120642
$ws.Range("I8").Value = 0.7652133098225272
$ws.Range("J8").Value = 0.6848199993188243
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.42828
$ws.Range("N8").Value = 0.85656
$ws.Range("O8").Value = 0.123391035029171
$ws.Range("P8").Value = 0.08578923067523865
$ws.Range("Q8").Value = 0.02583427788
$ws.Range("R8").Value = 0.10333711152
$ws.Range("S8").Value = 0.09442046231709933
$ws.Range("T8").Value = 0.05875018089257939

# Row 9: MuSCs -> Resolving-Mac (new row)
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Fgf15"
$ws.Range("C9").Value = "Fgfr3"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.060321
$ws.Range("H9").Value = 0.120642
$ws.Range("I9").Value = 0.7652133098225272
$ws.Range("J9").Value = 0.6848199993188243
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.05413033333333334
$ws.Range("N9").Value = 0.162391
$ws.Range("O9").Value = 0.01559539987035126
$ws.Range("P9").Value = 0.01626435854882633
$ws.Range("Q9").Value = 0.003265195837
$ws.Range("R9").Value = 0.019591175022
$ws.Range("S9").Value = 0.0119338075527973
$ws.Range("T9").Value = 0.01113815801032836
